$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select column D (mirrors the user selecting the whole column before deleting)
$colD = $ws.Range("D1:D1048576")
$colD.Select()

# Clear the contents of column D (removes the LEFT() helper formulas)
$colD.ClearContents()
